$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.561.48'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '3.752.78'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '613.82'
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.16'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = '3.753.21'
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  -2.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.166'
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.55'
$ws.Range("E11").Value = '  +3.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.480'
$ws.Range("E12").Value = '  -3.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.89'
$ws.Range("E13").Value = '  -2.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000253'
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("D15").Value = '4.366.72'
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").Value = '3.744.31'
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("D17").Value = '69.604.37'
$ws.Range("E17").Value = '  -0.32%  '
$ws.Range("E18").Value = '  -2.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.43'
$ws.Range("E19").Value = '  -2.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.37'
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '500.17'
$ws.Range("E21").Value = '  -3.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.17'
$ws.Range("E22").Value = '  -4.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.719'
$ws.Range("E23").Value = '  -1.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.58'
$ws.Range("E24").Value = '  +2.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.85'
$ws.Range("E25").Value = '  -2.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.88'
$ws.Range("E26").Value = '  -3.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.10'
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000134'
$ws.Range("E28").Value = '  +5.46%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.47'
$ws.Range("E30").Value = '  -1.71%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.08'
$ws.Range("E31").Value = '  +2.77%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.90'
$ws.Range("E32").Value = '  +2.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.38'
$ws.Range("E33").Value = '  -3.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.113'
$ws.Range("E34").Value = '  -2.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.04'
$ws.Range("E36").Value = '  +0.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.10'
$ws.Range("E37").Value = '  -1.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.348'
$ws.Range("E38").Value = '  +2.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.138'
$ws.Range("E39").Value = '  +3.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.06'
$ws.Range("E40").Value = '  +11.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '443.88'
$ws.Range("E41").Value = '  +5.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.06'
$ws.Range("E42").Value = '  -5.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.70'
$ws.Range("E43").Value = '  -2.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '44.46'
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.54'
$ws.Range("E45").Value = '  -3.33%  '
$ws.Range("D46").Value = '2.949.55'
$ws.Range("E46").Value = '  -4.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0358'
$ws.Range("E47").Value = '  -1.80%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.23'
$ws.Range("E49").Value = '  +1.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '27.03'
$ws.Range("E50").Value = '  -2.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.47'
$ws.Range("E51").Value = '  -2.22%  '
